$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.055.30"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.835.38"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "243.07"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "0.6273"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.07605"
$ws.Range("E8").Value = "  +3.70%  "
$ws.Range("D9").Value = "0.2932"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.60"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "1.851.85"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.960"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "0.6647"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001008"
$ws.Range("E15").Value = "  +17.16%  "
$ws.Range("D16").Value = "82.83"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "6.058"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "29.071.58"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "227.27"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").Value = "12.37"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "7.214"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.10"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "8.516"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "0.1384"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "17.93"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "1.498"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "4.103"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.020"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "1.193"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "0.05249"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "1.842"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "1.138"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").Value = "2.697"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").Value = "1.245.91"
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.760"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "0.01786"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").Value = "6.359"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "0.8977"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "102.17"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "1.982.47"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("D46").Value = "64.38"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5110"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").Value = "0.4043"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").Value = "8.884"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "0.05754"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "6.684"
$ws.Range("E51").Value = "  +0.15%  "
